$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Zoom level used while filling in the test plan
$ws.Application.ActiveWindow.Zoom = 72

# Row 7 - __init__ / Attributes set to input values.
$ws.Range("E7").Value = "none"
$ws.Range("F7").Value = '(7,"Beerdavinder", "singh","beerdavinder@pixel.com")'
$ws.Range("G7").Value = "Attributes are set"

# Fill in the rest of column E (Preconditions) for rows 8-11
$ws.Range("E8").Value = "none"
$ws.Range("E9").Value = "none"
$ws.Range("E10").Value = "none"
$ws.Range("E11").Value = "none"

# Fill in the rest of column F (Method Inputs) for rows 8-11
$ws.Range("F8").Value = '("bd","Beerdavinder", "singh","beerdavinder@pixel.com")'
$ws.Range("F9").Value = '(7,"", "singh","beerdavinder@pixel.com")'
$ws.Range("F10").Value = '(7,"Beerdavinder", "","beerdavinder@pixel.com")'
$ws.Range("F11").Value = '(7,"Beerdavinder", "singh","invalid")'

# Fill in the rest of column G (Expected Result) for rows 8-11
$ws.Range("G8").Value = "raise valueerror"
$ws.Range("G9").Value = "raise valueerror"
$ws.Range("G10").Value = "raise valueerror"
$ws.Range("G11").Value = "raise valueerror"

# Column G for rows 12-15 (Expected Result)
$ws.Range("G12").Value = "Return 7"
$ws.Range("G13").Value = "return beerdavinder"
$ws.Range("G14").Value = "return singh"
$ws.Range("G15").Value = "return beerdavinder@pixel.com"

# Column E for rows 12-16 (Preconditions)
$ws.Range("E12").Value = '(7,"Beerdavinder", "singh","beerdavinder@pixel.com")'
$ws.Range("E13").Value = '(7,"Beerdavinder", "singh","beerdavinder@pixel.com")'
$ws.Range("E14").Value = '(7,"Beerdavinder", "singh","beerdavinder@pixel.com")'
$ws.Range("E15").Value = '(7,"Beerdavinder", "singh","beerdavinder@pixel.com")'
$ws.Range("E16").Value = '(7,"Beerdavinder", "singh","beerdavinder@pixel.com")'

# Column F for rows 12-15 (Method Inputs)
$ws.Range("F12").Value = "self.client_number"
$ws.Range("F13").Value = "self.first_name"
$ws.Range("F14").Value = "self.last_name"
$ws.Range("F15").Value = "self.email_address"

# Row 16 - __str__ / Returns string in expected format.
$ws.Range("F16").Value = "__str__"
$ws.Range("G16").Value = 'Returns "Singh, Beerdavinder [7] - beerdavinder@pixel.com'

# F12 lost its wrap/bold formatting (content now fits without wrapping)
$ws.Range("F12").WrapText = $false
$ws.Range("F12").Font.Bold = $false

$ws.Range("F16").Select()
